# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11
$ws.Range("O11").Value = 1.29
$ws.Range("P11").Value = 3.5
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = 1.85

# Row 12
$ws.Range("G12").Value = 1.27
$ws.Range("I12").Value = 12
$ws.Range("L12").Value = 10
$ws.Range("U12").Value = 2.63
$ws.Range("V12").Value = 1.44
$ws.Range("Z12").Value = 7
$ws.Range("AC12").Value = 9.5
$ws.Range("AG12").Value = 21
$ws.Range("AI12").Value = 34
$ws.Range("AK12").Value = 101
$ws.Range("AN12").Value = 3
$ws.Range("AW12").Value = 11
$ws.Range("AZ12").Value = 351
$ws.Range("BA12").Value = 351

# Row 13
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 6.95
